# Final commit of example code
# Turns the blank "Sheet1" into the "In stock, payment accepted" spec-by-example sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "In stock, payment accepted"

# --- Given/When/Then narrative (rows 1-5) ---
$ws.Cells.Item(1,1).Value = "GIVEN"
$ws.Cells.Item(1,2).Value = "a CD that's in the Top 100 and we have it in stock, and the customer's card payment will be accepted"

$ws.Cells.Item(2,1).Value = "WHEN"
$ws.Cells.Item(2,2).Value = "The customer buy that CD"

$ws.Cells.Item(3,1).Value = "THEN"
$ws.Cells.Item(3,2).Value = "One copy is deducted from CD's stock"

$ws.Cells.Item(4,1).Value = "AND"
$ws.Cells.Item(4,2).Value = "The customer's card is charged £1 less than the lowest price from competitors found for that CD"

$ws.Cells.Item(5,1).Value = "AND"
$ws.Cells.Item(5,2).Value = "The charts are notified of the sale"

# --- Inputs / Outputs labels (row 7), bold (bold font is registered first so it
# ends up as fontId 1, matching the authored file) ---
$ws.Cells.Item(7,1).Value = "Inputs"
$ws.Cells.Item(7,8).Value = "Outputs"
$ws.Range("A7").Font.Bold = $true
$ws.Range("H7").Font.Bold = $true

# Column A of rows 1-5 is italic
$ws.Range("A1:A5").Font.Italic = $true

# --- Table header row (row 8) ---
$ws.Cells.Item(8,1).Value = "CD"
$ws.Cells.Item(8,2).Value = "Chart Position"
$ws.Cells.Item(8,3).Value = "Stock"
$ws.Cells.Item(8,4).Value = "Price"
$ws.Cells.Item(8,5).Value = "Credit Card"
$ws.Cells.Item(8,6).Value = "Payment Accepted"
$ws.Cells.Item(8,7).Value = "Lowest Price"
$ws.Cells.Item(8,8).Value = "End Stock"
$ws.Cells.Item(8,9).Value = "Charged"
$ws.Cells.Item(8,10).Value = "Notification"

# --- Table data row (row 9) ---
$ws.Cells.Item(9,1).Value = 'title: "NumberOf The Beast", artist: "Iron Maiden"'
$ws.Cells.Item(9,2).Value = 100
$ws.Cells.Item(9,3).Value = 10
$ws.Cells.Item(9,4).Value = 9.99
$ws.Cells.Item(9,5).Value = 1234234634564560
$ws.Cells.Item(9,6).Value = $true
$ws.Cells.Item(9,7).Value = 7.99
$ws.Cells.Item(9,8).Value = 9
$ws.Cells.Item(9,9).Value = 6.99
$ws.Cells.Item(9,10).Value = 'sales: 1, album: "Iron Maiden - Number Of The Beast"'

# --- Turn A8:J9 into an Excel Table ---
$tbl = $ws.ListObjects.Add(1, $ws.Range("A8:J9"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = "TableStyleLight8"

# --- Column widths (character widths; engine adds the standard 5/6 cell padding) ---
$ws.Columns.Item(1).ColumnWidth = 44.053385416666664
$ws.Columns.Item(2).ColumnWidth = 20.721354166666668
$ws.Columns.Item(3).ColumnWidth = 11.830729166666666
$ws.Columns.Item(4).ColumnWidth = 11.830729166666666
$ws.Columns.Item(5).ColumnWidth = 19.608072916666668
$ws.Columns.Item(6).ColumnWidth = 20.385416666666668
$ws.Columns.Item(7).ColumnWidth = 14.721354166666666
$ws.Columns.Item(8).ColumnWidth = 12.608072916666666
$ws.Columns.Item(9).ColumnWidth = 13.385416666666666
$ws.Columns.Item(10).ColumnWidth = 11.830729166666666

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection shown in the saved view ---
$ws.Range("A8:J9").Select() | Out-Null
